# Master Data Tables - Test Data: master-reg_center_device_h.xlsx
# Append 45 new rows (rows 102-146) to the "master-reg_center_device_h" sheet,
# continuing the existing regcntr_id/device_id/lang_code/is_active/cr_by/
# cr_dtimes/eff_dtimes pattern, then update the view to match where the
# author left the selection/scroll position, and set the page orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 102
$rowCount = 45
$startDeviceId = 3000121

# regcntr_id cycles through this 9-value pattern for the new block of rows
$regCenterPattern = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)

$data = New-Object 'object[,]' $rowCount,7
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = $regCenterPattern[$i % $regCenterPattern.Length]
    $data[$i,1] = $startDeviceId + $i
    $data[$i,2] = "eng"
    $data[$i,3] = $true
    $data[$i,4] = "superadmin"
    $data[$i,5] = "now()"
    $data[$i,6] = "now()"
}

$endRow = $startRow + $rowCount - 1
$targetRange = $ws.Range("A" + $startRow + ":G" + $endRow)
$targetRange.Value = $data

# Reflect the author's final selection/scroll position on the sheet.
$selRange = $ws.Range("A" + $startRow + ":B" + $endRow)
[void]$selRange.Select()
$excel.ActiveWindow.ScrollRow = 128
$excel.ActiveWindow.ScrollColumn = 1

# Page setup: orientation was switched to portrait.
$ws.PageSetup.Orientation = 1
